# Update the Cryptos worksheet with the latest scraped price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (matches the sheet's existing
# inline-string cells) instead of letting Excel auto-coerce numeric-looking
# strings (e.g. "1.0000", "0.2792") into numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "30.437.47"
Set-TextValue $ws.Range("E2") "  +0.48%  "
Set-TextValue $ws.Range("D3") "1.866.31"
Set-TextValue $ws.Range("E3") "  -0.11%  "
Set-TextValue $ws.Range("E4") "  -0.06%  "
Set-TextValue $ws.Range("D5") "235.69"
Set-TextValue $ws.Range("D6") "1.0000"
Set-TextValue $ws.Range("E6") "  -0.11%  "
Set-TextValue $ws.Range("D7") "0.4787"
Set-TextValue $ws.Range("E7") "  -0.04%  "
Set-TextValue $ws.Range("D8") "0.2792"
Set-TextValue $ws.Range("D9") "0.06542"
Set-TextValue $ws.Range("E9") "  +0.73%  "
Set-TextValue $ws.Range("D10") "1.856.03"
Set-TextValue $ws.Range("E10") "  -0.56%  "
Set-TextValue $ws.Range("D11") "0.07443"
Set-TextValue $ws.Range("E11") "  +0.06%  "
Set-TextValue $ws.Range("D12") "16.20"
Set-TextValue $ws.Range("E12") "  -1.80%  "
Set-TextValue $ws.Range("D13") "5.066"
Set-TextValue $ws.Range("E13") "  -0.03%  "
Set-TextValue $ws.Range("D14") "86.87"
Set-TextValue $ws.Range("E14") "  -1.22%  "
Set-TextValue $ws.Range("D15") "0.6389"
Set-TextValue $ws.Range("E15") "  -2.17%  "
Set-TextValue $ws.Range("D16") "30.423.09"
Set-TextValue $ws.Range("E16") "  +0.53%  "
Set-TextValue $ws.Range("D17") "1.000"
Set-TextValue $ws.Range("E17") "  -0.10%  "
Set-TextValue $ws.Range("D18") "12.97"
Set-TextValue $ws.Range("E18") "  -2.23%  "
Set-TextValue $ws.Range("D19") "232.25"
Set-TextValue $ws.Range("E19") "  +5.88%  "
Set-TextValue $ws.Range("D20") "0.000007456"
Set-TextValue $ws.Range("E20") "  -1.28%  "
Set-TextValue $ws.Range("D21") "2.118.97"
Set-TextValue $ws.Range("E21") "  +0.59%  "
Set-TextValue $ws.Range("E22") "  -0.05%  "
Set-TextValue $ws.Range("E23") "  -2.78%  "
Set-TextValue $ws.Range("D24") "6.081"
Set-TextValue $ws.Range("E24") "  -1.09%  "
Set-TextValue $ws.Range("D25") "168.37"
Set-TextValue $ws.Range("E25") "  +0.01%  "
Set-TextValue $ws.Range("D26") "9.302"
Set-TextValue $ws.Range("E26") "  -0.01%  "
Set-TextValue $ws.Range("D27") "18.14"
Set-TextValue $ws.Range("E27") "  -1.58%  "
Set-TextValue $ws.Range("D28") "1.897"
Set-TextValue $ws.Range("E28") "  -3.58%  "
Set-TextValue $ws.Range("D29") "0.1044"
Set-TextValue $ws.Range("E29") "  +11.42%  "
Set-TextValue $ws.Range("D30") "1.381"
Set-TextValue $ws.Range("E30") "  -5.09%  "
Set-TextValue $ws.Range("D31") "4.258"
Set-TextValue $ws.Range("E31") "  -0.91%  "
Set-TextValue $ws.Range("D32") "3.966"
Set-TextValue $ws.Range("E32") "  -1.10%  "
Set-TextValue $ws.Range("D33") "0.04970"
Set-TextValue $ws.Range("E33") "  -1.59%  "
Set-TextValue $ws.Range("D34") "1.170"
Set-TextValue $ws.Range("E34") "  -2.70%  "
Set-TextValue $ws.Range("D35") "0.7399"
Set-TextValue $ws.Range("E35") "  -1.06%  "
Set-TextValue $ws.Range("D36") "0.9996"
Set-TextValue $ws.Range("E37") "  -0.02%  "
Set-TextValue $ws.Range("D38") "0.01942"
Set-TextValue $ws.Range("E38") "  +6.77%  "
Set-TextValue $ws.Range("D39") "2.637"
Set-TextValue $ws.Range("E39") "  +1.03%  "
Set-TextValue $ws.Range("D40") "0.9133"
Set-TextValue $ws.Range("E40") "  +0.96%  "
Set-TextValue $ws.Range("D41") "2.034"
Set-TextValue $ws.Range("E41") "  -2.01%  "
Set-TextValue $ws.Range("D42") "106.12"
Set-TextValue $ws.Range("E42") "  -0.53%  "
Set-TextValue $ws.Range("D43") "0.9961"
Set-TextValue $ws.Range("E43") "  -0.69%  "
Set-TextValue $ws.Range("D44") "0.4173"
Set-TextValue $ws.Range("E44") "  -2.03%  "
Set-TextValue $ws.Range("D45") "5.565"
Set-TextValue $ws.Range("E45") "  -6.39%  "
Set-TextValue $ws.Range("D46") "7.157"
Set-TextValue $ws.Range("E46") "  -2.87%  "
Set-TextValue $ws.Range("D47") "61.77"
Set-TextValue $ws.Range("E47") "  -3.10%  "
Set-TextValue $ws.Range("E48") "  -4.52%  "
Set-TextValue $ws.Range("D49") "8.884"
Set-TextValue $ws.Range("E49") "  -1.16%  "
Set-TextValue $ws.Range("D50") "1.411"
Set-TextValue $ws.Range("E50") "  -3.92%  "
Set-TextValue $ws.Range("D51") "33.42"
Set-TextValue $ws.Range("E51") "  -0.47%  "
